# edit.ps1 - reproduce the author's edit:
#   1) Footer "last updated" date placeholders bump from 2018/11/13 -> 2018/11/19
#      (present on the slide master and every slide layout).
#   2) On the "Current branch" slide (slide 16), three adjacent runs that were
#      split mid-sentence are re-merged into single runs (no visible text
#      change, just consolidating the run structure).
#   3) On the "switch/upload branch" slide (slide 17), three adjacent runs of
#      the footnote are likewise re-merged into a single run.
#
# Helper: given a TextRange and a target substring, relocate that exact
# substring inside the range and re-assign its Characters(...).Text to itself.
# Re-assigning a sub-range forces the engine to rebuild that span as a single
# run (inheriting the formatting of the first run touched), which is exactly
# how adjacent runs with identical formatting get coalesced into one.
function Merge-Span {
    param(
        [object]$Range,
        [string]$Needle
    )
    $full = $Range.Text
    $idx = $full.IndexOf($Needle)
    if ($idx -ge 0) {
        $chars = $Range.Characters($idx + 1, $Needle.Length)
        $chars.Text = $Needle
    }
}

$p = $ppt.ActivePresentation

# --- 1) Update the cached "datetimeFigureOut" footer date everywhere it
#        appears: the slide master plus every custom (slide) layout. ---
$master = $p.SlideMaster

$dateContainers = New-Object System.Collections.ArrayList
[void]$dateContainers.Add($master)
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    [void]$dateContainers.Add($master.CustomLayouts.Item($li))
}

foreach ($container in $dateContainers) {
    for ($si = 1; $si -le $container.Shapes.Count; $si++) {
        $shp = $container.Shapes.Item($si)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $t = $tr.Text
            if ($t -like "*2018/11/13*") {
                $tr.Text = $t.Replace("2018/11/13", "2018/11/19")
            }
        }
    }
}

# --- 2) Slide 16 ("切换分支完成后点击..."): merge three run pairs back
#        together. ---
$slide16 = $p.Slides.Item(16)
$shape16 = $slide16.Shapes.Item(2)
$tr16 = $shape16.TextFrame.TextRange

Merge-Span $tr16 "Current branch-New branch"
Merge-Span $tr16 "点击菜单栏的"
Merge-Span $tr16 "），输入分支名，点击"

# --- 3) Slide 17 ("注：上传分支中的单个文件..."): merge three runs into
#        one. ---
$slide17 = $p.Slides.Item(17)
$shape17 = $slide17.Shapes.Item(2)
$tr17 = $shape17.TextFrame.TextRange

Merge-Span $tr17 "注：上传分支中的单个文件大小一般不能超过"
